$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 683.8
$ws.Range("I2").Value = 683.8
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 683.8
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -570.8
$ws.Range("H17").Value = 308358.6
$ws.Range("J17").Value = 313862.1
$ws.Range("L17").Value = 941586.2999999999
$ws.Range("N17").Value = -941922.2999999999
$ws.Range("H40").Value = 1771
$ws.Range("I40").Value = 1784.4445
$ws.Range("J40").Value = 1650
$ws.Range("K40").Value = 1784.4445
$ws.Range("L40").Value = 1650
$ws.Range("M40").Value = -1609.4445
$ws.Range("N40").Value = -2000
$ws.Range("H121").Value = 605.0323
$ws.Range("J121").Value = 595.37933
$ws.Range("L121").Value = 1786.13799
$ws.Range("N121").Value = -5280.13799
$ws.Range("H129").Value = 939.4643
$ws.Range("I129").Value = 493.33334
$ws.Range("J129").Value = 993
$ws.Range("K129").Value = 1480.00002
$ws.Range("L129").Value = 2979
$ws.Range("M129").Value = 3519.99998
$ws.Range("N129").Value = -12979
$ws.Range("H137").Value = 2453.125
$ws.Range("I137").Value = 1305.5555
$ws.Range("K137").Value = 3916.6665
$ws.Range("M137").Value = -1366.6665
$ws.Range("H138").Value = 2220.1333
$ws.Range("I138").Value = 930
$ws.Range("J138").Value = 3510.2666
$ws.Range("K138").Value = 2790
$ws.Range("L138").Value = 10530.7998
$ws.Range("M138").Value = 2350
$ws.Range("N138").Value = -20810.7998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10505.9375
$ws.Range("I45").Value = 15943.6
$ws.Range("J45").Value = 1443.1666
$ws.Range("K45").Value = 15943.6
$ws.Range("L45").Value = 1443.1666
$ws.Range("M45").Value = -15566.6
$ws.Range("N45").Value = -2197.1666
$ws.Range("H63").Value = 76925140
$ws.Range("I63").Value = 111112890
$ws.Range("J63").Value = 2700
$ws.Range("K63").Value = 111112890
$ws.Range("L63").Value = 2700
$ws.Range("M63").Value = -111112204
$ws.Range("N63").Value = -4072
$ws.Range("H66").Value = 76925140
$ws.Range("I66").Value = 111112890
$ws.Range("J66").Value = 2700
$ws.Range("K66").Value = 555564450
$ws.Range("L66").Value = 13500
$ws.Range("M66").Value = -555561018
$ws.Range("N66").Value = -20364
$ws.Range("H110").Value = 1066.3636
$ws.Range("I110").Value = 953
$ws.Range("J110").Value = 1309.2858
$ws.Range("K110").Value = 953
$ws.Range("L110").Value = 1309.2858
$ws.Range("M110").Value = 1092
$ws.Range("N110").Value = -5399.2858
$ws.Range("H132").Value = 2795.4
$ws.Range("I132").Value = 1748.2727
$ws.Range("J132").Value = 4075.2222
$ws.Range("K132").Value = 5244.8181
$ws.Range("L132").Value = 12225.6666
$ws.Range("M132").Value = -2714.8181
$ws.Range("N132").Value = -17285.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1156.5227
$ws.Range("I94").Value = 932.63635
$ws.Range("K94").Value = 932.63635
$ws.Range("M94").Value = -481.63635

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2055.823
$ws.Range("I31").Value = 1431.3
$ws.Range("J31").Value = 2339.697
$ws.Range("K31").Value = 1431.3
$ws.Range("L31").Value = 2339.697
$ws.Range("M31").Value = -1136.3
$ws.Range("N31").Value = -2929.697
$ws.Range("H34").Value = 2055.823
$ws.Range("I34").Value = 1431.3
$ws.Range("J34").Value = 2339.697
$ws.Range("K34").Value = 1431.3
$ws.Range("L34").Value = 2339.697
$ws.Range("M34").Value = -1229.3
$ws.Range("N34").Value = -2743.697
$ws.Range("H86").Value = 83335530
$ws.Range("I86").Value = 125002390
$ws.Range("J86").Value = 1794.25
$ws.Range("K86").Value = 125002390
$ws.Range("L86").Value = 1794.25
$ws.Range("M86").Value = -125001267
$ws.Range("N86").Value = -4040.25
$ws.Range("H89").Value = 83335530
$ws.Range("I89").Value = 125002390
$ws.Range("J89").Value = 1794.25
$ws.Range("K89").Value = 625011950
$ws.Range("L89").Value = 8971.25
$ws.Range("M89").Value = -625006334
$ws.Range("N89").Value = -20203.25
$ws.Range("H99").Value = 7365068
$ws.Range("I99").Value = 11076
$ws.Range("J99").Value = 17870772
$ws.Range("K99").Value = 11076
$ws.Range("L99").Value = 17870772
$ws.Range("M99").Value = -9578
$ws.Range("N99").Value = -17873768
$ws.Range("H125").Value = 48881.25
$ws.Range("J125").Value = 48881.25
$ws.Range("L125").Value = 48881.25
$ws.Range("N125").Value = -53801.25
$ws.Range("H126").Value = 7365068
$ws.Range("I126").Value = 11076
$ws.Range("J126").Value = 17870772
$ws.Range("K126").Value = 33228
$ws.Range("L126").Value = 53612316
$ws.Range("M126").Value = -30758
$ws.Range("N126").Value = -53617256
$ws.Range("H132").Value = 1477.5385
$ws.Range("I132").Value = 1199.9395
$ws.Range("J132").Value = 3004.3333
$ws.Range("K132").Value = 3599.8185
$ws.Range("L132").Value = 9012.999899999999
$ws.Range("M132").Value = -1069.8185
$ws.Range("N132").Value = -14072.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3564.8
$ws.Range("J68").Value = 2449.9167
$ws.Range("L68").Value = 7349.750100000001
$ws.Range("N68").Value = -8971.750100000001
$ws.Range("H71").Value = 3564.8
$ws.Range("J71").Value = 2449.9167
$ws.Range("L71").Value = 22049.2503
$ws.Range("N71").Value = -30161.2503
$ws.Range("H76").Value = 5000
$ws.Range("J76").Value = 5000
$ws.Range("L76").Value = 15000
$ws.Range("N76").Value = -15766
$ws.Range("H79").Value = 5000
$ws.Range("J79").Value = 5000
$ws.Range("L79").Value = 15000
$ws.Range("N79").Value = -17652
$ws.Range("H107").Value = 776.7455
$ws.Range("J107").Value = 1059.1428
$ws.Range("L107").Value = 3177.4284
$ws.Range("N107").Value = -7017.428400000001
$ws.Range("H122").Value = 490.375
$ws.Range("I122").Value = 476.5
$ws.Range("J122").Value = 495
$ws.Range("K122").Value = 4288.5
$ws.Range("L122").Value = 4455
$ws.Range("M122").Value = -1838.5
$ws.Range("N122").Value = -9355

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 25499.5
$ws.Range("J96").Value = 25499.5
$ws.Range("L96").Value = 25499.5
$ws.Range("N96").Value = -30991.5
$ws.Range("H97").Value = 755.4545000000001
$ws.Range("I97").Value = 733.3333
$ws.Range("J97").Value = 782
$ws.Range("K97").Value = 733.3333
$ws.Range("L97").Value = 782
$ws.Range("M97").Value = -237.3333
$ws.Range("N97").Value = -1774

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 10571653
$ws.Range("I132").Value = 14944817
$ws.Range("J132").Value = 3174.4167
$ws.Range("K132").Value = 44834451
$ws.Range("L132").Value = 9523.250100000001
$ws.Range("M132").Value = -44831921
$ws.Range("N132").Value = -14583.2501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4982.3335
$ws.Range("I122").Value = 6021.8887
$ws.Range("J122").Value = 3942.7778
$ws.Range("K122").Value = 18065.6661
$ws.Range("L122").Value = 11828.3334
$ws.Range("M122").Value = -15615.6661
$ws.Range("N122").Value = -16728.3334
